$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.024.09'
$ws.Range("E2").Value = '  -0.11%  '

$ws.Range("D3").Value = '''1.833.07'
$ws.Range("E3").Value = '  +0.22%  '

$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = '''244.58'
$ws.Range("E5").Value = '  +1.54%  '

$ws.Range("D6").Value = '''0.6339'
$ws.Range("E6").Value = '  +1.87%  '

$ws.Range("D7").Value = '''1.0000'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '''0.07582'
$ws.Range("E8").Value = '  +2.95%  '

$ws.Range("D9").Value = '''0.2949'
$ws.Range("E9").Value = '  +0.90%  '

$ws.Range("D10").Value = '''22.80'
$ws.Range("E10").Value = '  +0.46%  '

$ws.Range("D11").Value = '''0.07740'
$ws.Range("E11").Value = '  +1.25%  '

$ws.Range("D12").Value = '''1.828.17'
$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("D13").Value = '''4.998'
$ws.Range("E13").Value = '  +0.68%  '

$ws.Range("D14").Value = '''0.6709'
$ws.Range("E14").Value = '  +1.18%  '

$ws.Range("D15").Value = '''83.19'
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("D16").Value = '''0.000009832'
$ws.Range("E16").Value = '  +8.85%  '

$ws.Range("D17").Value = '''6.119'
$ws.Range("E17").Value = '  +1.46%  '

$ws.Range("D18").Value = '''29.040.67'
$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("E19").Value = '  +1.58%  '

$ws.Range("D20").Value = '''226.62'
$ws.Range("E20").Value = '  +0.52%  '

$ws.Range("D22").Value = '''7.223'
$ws.Range("E22").Value = '  +0.86%  '

$ws.Range("D23").Value = '''1.001'
$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").Value = '''160.30'
$ws.Range("E24").Value = '  +0.44%  '

$ws.Range("D25").Value = '''0.1405'
$ws.Range("E25").Value = '  +3.51%  '

$ws.Range("D26").Value = '''8.537'
$ws.Range("E26").Value = '  +1.33%  '

$ws.Range("D27").Value = '''17.95'
$ws.Range("E27").Value = '  +0.80%  '

$ws.Range("E28").Value = '  +0.34%  '

$ws.Range("D29").Value = '''4.124'
$ws.Range("E29").Value = '  +1.66%  '

$ws.Range("D30").Value = '''4.057'
$ws.Range("E30").Value = '  +0.52%  '

$ws.Range("D32").Value = '''0.05393'
$ws.Range("E32").Value = '  +2.73%  '

$ws.Range("D33").Value = '''1.860'
$ws.Range("E33").Value = '  +0.88%  '

$ws.Range("D34").Value = '''0.7474'
$ws.Range("E34").Value = '  +1.96%  '

$ws.Range("E35").Value = '  -0.96%  '

$ws.Range("D36").Value = '''2.666'
$ws.Range("E36").Value = '  +0.75%  '

$ws.Range("D37").Value = '''1.237.94'
$ws.Range("E37").Value = '  -4.14%  '

$ws.Range("E38").Value = '  +0.56%  '

$ws.Range("D39").Value = '''2.760'
$ws.Range("E39").Value = '  +0.50%  '

$ws.Range("D40").Value = '''6.628'
$ws.Range("E40").Value = '  +4.99%  '

$ws.Range("D41").Value = '''0.9030'
$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("D42").Value = '''0.9996'
$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("D43").Value = '''102.60'
$ws.Range("E43").Value = '  +0.71%  '

$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '''1.984.28'
$ws.Range("E44").Value = '  +0.40%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '''0.00000000125'
$ws.Range("E45").Value = '  +4.69%  '

$ws.Range("D46").Value = '''64.87'
$ws.Range("E46").Value = '  +1.25%  '

$ws.Range("D47").Value = '''0.5113'
$ws.Range("E47").Value = '  -0.03%  '

$ws.Range("D48").Value = '''0.4097'
$ws.Range("E48").Value = '  +3.32%  '

$ws.Range("D49").Value = '''9.065'
$ws.Range("E49").Value = '  +2.47%  '

$ws.Range("D50").Value = '''6.777'
$ws.Range("E50").Value = '  +1.70%  '

$ws.Range("D51").Value = '''0.05783'
$ws.Range("E51").Value = '  +0.28%  '
